# Initial Data File Update
# Adds 13 new transaction rows (160-172) to the "Transacciones" sheet,
# mirroring the bank/expense ledger pattern already present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transacciones")

# ---------------------------------------------------------------------
# 1) Create rows 160-172 by inserting a copy of the row immediately
#    above each new row. This clones number formats / styles (date
#    format on column A, the "default" styles on N/O, the "Good" style
#    on P) without minting brand-new style entries, and keeps the
#    formulas' relative references correctly anchored to their own row.
# ---------------------------------------------------------------------
for ($i = 160; $i -le 172; $i++) {
    $ws.Rows.Item($i - 1).Copy()
    $ws.Rows.Item($i).Insert(-4121)
}

# ---------------------------------------------------------------------
# 2) Overwrite each new row with the real transaction data.
# ---------------------------------------------------------------------

# Row 160
$ws.Range("A160").Value2 = 43595
$ws.Range("B160").Value2 = 18.9
$ws.Range("C160").Value2 = "Almuerzo Foraneo"
$ws.Range("D160").Value2 = "Comida"
$ws.Range("E160").Value2 = "Gasto"
$ws.Range("F160").Value2 = "Tarjeta Santander"
$ws.Range("G160").Value2 = "Extra"
$ws.Range("K160").Value2 = 5064.18
$ws.Range("L160").Formula = "=L159-B160"
$ws.Range("M160").Value2 = 5
$ws.Range("N160").Formula = "=SUM(K160:M160)"
$ws.Range("O160").Formula = "=N160-4000"
$ws.Range("P160").Formula = "=O160-Ahorros!`$E`$4"

# Row 161
$ws.Range("A161").Value2 = 43595
$ws.Range("B161").Value2 = 223.96
$ws.Range("C161").Value2 = "Gasolina"
$ws.Range("D161").Value2 = "Gasolina"
$ws.Range("E161").Value2 = "Gasto"
$ws.Range("F161").Value2 = "Tarjeta Santander"
$ws.Range("G161").Value2 = "Gasolinería Mobil"
$ws.Range("K161").Value2 = 5064.18
$ws.Range("L161").Formula = "=L160-B161"
$ws.Range("M161").Value2 = 5
$ws.Range("N161").Formula = "=SUM(K161:M161)"
$ws.Range("O161").Formula = "=N161-4000"
$ws.Range("P161").Formula = "=O161-Ahorros!`$E`$4"

# Row 162
$ws.Range("A162").Value2 = 43596
$ws.Range("B162").Value2 = 32.5
$ws.Range("C162").Value2 = "Paleta Magnum"
$ws.Range("D162").Value2 = "Golosinas"
$ws.Range("E162").Value2 = "Gasto"
$ws.Range("F162").Value2 = "Tarjeta Banamex"
$ws.Range("G162").Value2 = "Oxxo"
$ws.Range("K162").Formula = "=K161-B162"
$ws.Range("L162").Value2 = 3196.26
$ws.Range("M162").Value2 = 5
$ws.Range("N162").Formula = "=SUM(K162:M162)"
$ws.Range("O162").Formula = "=N162-4000"
$ws.Range("P162").Formula = "=O162-Ahorros!`$E`$4"

# Row 163
$ws.Range("A163").Value2 = 43596
$ws.Range("B163").Value2 = 26.5
$ws.Range("C163").Value2 = "Refrescos "
$ws.Range("D163").Value2 = "Golosinas"
$ws.Range("E163").Value2 = "Gasto"
$ws.Range("F163").Value2 = "Tarjeta Banamex"
$ws.Range("G163").Value2 = "Farmacia Guadalajara"
$ws.Range("K163").Formula = "=K162-B163"
$ws.Range("L163").Value2 = 3196.26
$ws.Range("M163").Value2 = 5
$ws.Range("N163").Formula = "=SUM(K163:M163)"
$ws.Range("O163").Formula = "=N163-4000"
$ws.Range("P163").Formula = "=O163-Ahorros!`$E`$4"

# Row 164
$ws.Range("A164").Value2 = 43596
$ws.Range("B164").Value2 = 149
$ws.Range("C164").Value2 = "Tarjeta Micro SD 16 GB"
$ws.Range("D164").Value2 = "Electrónicos"
$ws.Range("E164").Value2 = "Gasto"
$ws.Range("F164").Value2 = "Tarjeta Banamex"
$ws.Range("G164").Value2 = "Coppel"
$ws.Range("K164").Formula = "=K163-B164"
$ws.Range("L164").Value2 = 3196.26
$ws.Range("M164").Value2 = 5
$ws.Range("N164").Formula = "=SUM(K164:M164)"
$ws.Range("O164").Formula = "=N164-4000"
$ws.Range("P164").Formula = "=O164-Ahorros!`$E`$4"

# Row 165
$ws.Range("A165").Value2 = 43598
$ws.Range("B165").Value2 = 25.9
$ws.Range("C165").Value2 = "Almuerzo Foraneo"
$ws.Range("D165").Value2 = "Comida"
$ws.Range("E165").Value2 = "Gasto"
$ws.Range("F165").Value2 = "Tarjeta Santander"
$ws.Range("G165").Value2 = "Extra"
$ws.Range("K165").Value2 = 4856.18
$ws.Range("L165").Formula = "=L164-B165"
$ws.Range("M165").Value2 = 5
$ws.Range("N165").Formula = "=SUM(K165:M165)"
$ws.Range("O165").Formula = "=N165-4000"
$ws.Range("P165").Formula = "=O165-Ahorros!`$E`$4"

# Row 166
$ws.Range("A166").Value2 = 43598
$ws.Range("B166").Value2 = 12.5
$ws.Range("C166").Value2 = "Agua Mineral"
$ws.Range("D166").Value2 = "Despensa"
$ws.Range("E166").Value2 = "Gasto"
$ws.Range("F166").Value2 = "Tarjeta Santander"
$ws.Range("G166").Value2 = "Soriana"
$ws.Range("K166").Value2 = 4856.18
$ws.Range("L166").Formula = "=L165-B166"
$ws.Range("M166").Value2 = 5
$ws.Range("N166").Formula = "=SUM(K166:M166)"
$ws.Range("O166").Formula = "=N166-4000"
$ws.Range("P166").Formula = "=O166-Ahorros!`$E`$4"

# Row 167
$ws.Range("A167").Value2 = 43598
$ws.Range("B167").Value2 = 22.76
$ws.Range("C167").Value2 = "Chuleta de Cerdo"
$ws.Range("D167").Value2 = "Despensa"
$ws.Range("E167").Value2 = "Gasto"
$ws.Range("F167").Value2 = "Tarjeta Santander"
$ws.Range("G167").Value2 = "Soriana"
$ws.Range("K167").Value2 = 4856.18
$ws.Range("L167").Formula = "=L166-B167"
$ws.Range("M167").Value2 = 5
$ws.Range("N167").Formula = "=SUM(K167:M167)"
$ws.Range("O167").Formula = "=N167-4000"
$ws.Range("P167").Formula = "=O167-Ahorros!`$E`$4"

# Row 168
$ws.Range("A168").Value2 = 43598
$ws.Range("B168").Value2 = 38.21
$ws.Range("C168").Value2 = "Carne de Res"
$ws.Range("D168").Value2 = "Despensa"
$ws.Range("E168").Value2 = "Gasto"
$ws.Range("F168").Value2 = "Tarjeta Santander"
$ws.Range("G168").Value2 = "Soriana"
$ws.Range("K168").Value2 = 4856.18
$ws.Range("L168").Formula = "=L167-B168"
$ws.Range("M168").Value2 = 5
$ws.Range("N168").Formula = "=SUM(K168:M168)"
$ws.Range("O168").Formula = "=N168-4000"
$ws.Range("P168").Formula = "=O168-Ahorros!`$E`$4"

# Row 169
$ws.Range("A169").Value2 = 43598
$ws.Range("B169").Value2 = 12
$ws.Range("C169").Value2 = "Galletas Gamesa"
$ws.Range("D169").Value2 = "Despensa"
$ws.Range("E169").Value2 = "Gasto"
$ws.Range("F169").Value2 = "Tarjeta Santander"
$ws.Range("G169").Value2 = "Soriana"
$ws.Range("K169").Value2 = 4856.18
$ws.Range("L169").Formula = "=L168-B169"
$ws.Range("M169").Value2 = 5
$ws.Range("N169").Formula = "=SUM(K169:M169)"
$ws.Range("O169").Formula = "=N169-4000"
$ws.Range("P169").Formula = "=O169-Ahorros!`$E`$4"

# Row 170
$ws.Range("A170").Value2 = 43598
$ws.Range("B170").Value2 = 22
$ws.Range("C170").Value2 = "Leche Santa Clara"
$ws.Range("D170").Value2 = "Despensa"
$ws.Range("E170").Value2 = "Gasto"
$ws.Range("F170").Value2 = "Tarjeta Santander"
$ws.Range("G170").Value2 = "Soriana"
$ws.Range("K170").Value2 = 4856.18
$ws.Range("L170").Formula = "=L169-B170"
$ws.Range("M170").Value2 = 5
$ws.Range("N170").Formula = "=SUM(K170:M170)"
$ws.Range("O170").Formula = "=N170-4000"
$ws.Range("P170").Formula = "=O170-Ahorros!`$E`$4"

# Row 171
$ws.Range("A171").Value2 = 43598
$ws.Range("B171").Value2 = 15.16
$ws.Range("C171").Value2 = "Melón Chino Especial"
$ws.Range("D171").Value2 = "Despensa"
$ws.Range("E171").Value2 = "Gasto"
$ws.Range("F171").Value2 = "Tarjeta Santander"
$ws.Range("G171").Value2 = "Soriana"
$ws.Range("K171").Value2 = 4856.18
$ws.Range("L171").Formula = "=L170-B171"
$ws.Range("M171").Value2 = 5
$ws.Range("N171").Formula = "=SUM(K171:M171)"
$ws.Range("O171").Formula = "=N171-4000"
$ws.Range("P171").Formula = "=O171-Ahorros!`$E`$4"

# Row 172
$ws.Range("A172").Value2 = 43598
$ws.Range("B172").Value2 = 36.25
$ws.Range("C172").Value2 = "Pan Bimbo Integral"
$ws.Range("D172").Value2 = "Despensa"
$ws.Range("E172").Value2 = "Gasto"
$ws.Range("F172").Value2 = "Tarjeta Santander"
$ws.Range("G172").Value2 = "Soriana"
$ws.Range("K172").Value2 = 4856.18
$ws.Range("L172").Formula = "=L171-B172"
$ws.Range("M172").Value2 = 5
$ws.Range("N172").Formula = "=SUM(K172:M172)"
$ws.Range("O172").Formula = "=N172-4000"
$ws.Range("P172").Formula = "=O172-Ahorros!`$E`$4"

# ---------------------------------------------------------------------
# 3) Restore the sheet's scroll/selection state to reflect the new
#    bottom of the data (frozen header rows 1-3, scrolled down, with
#    the last active selection on P125).
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("A4").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("P125").Select()

$wb.Save()
